$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(86).Insert()

$ws.Cells.Item(86, 1).Value = 5
$ws.Cells.Item(86, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(86, 3).Value = "Maule"
$ws.Cells.Item(86, 4).Value = 44546
$ws.Cells.Item(86, 5).Value = 7
$ws.Cells.Item(86, 6).Value = "Fruta"
$ws.Cells.Item(86, 7).Value = 100108
$ws.Cells.Item(86, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(86, 9).Value = 100108005
$ws.Cells.Item(86, 10).Value = "Piña"
$ws.Cells.Item(86, 11).Value = "Caramelo"
$ws.Cells.Item(86, 12).Value = "Segunda"
$ws.Cells.Item(86, 13).Value = 300
$ws.Cells.Item(86, 14).Value = 17000
$ws.Cells.Item(86, 15).Value = 17000
$ws.Cells.Item(86, 16).Value = 17000
$ws.Cells.Item(86, 17).Value = '$/caja 14 unidades'
$ws.Cells.Item(86, 18).Value = "Ecuador"
$ws.Cells.Item(86, 19).Value = 1214
$ws.Cells.Item(86, 20).Value = 14
